# Auto-generated update of market-price columns (H..N) in the Coeurl_Profits
# workbook. Source data refreshed via scheduled market-data runner; each
# touched cell is written with its new literal value to match the refreshed
# snapshot (no formulas are involved -- these columns store plain numbers).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 210.58333
$ws.Range("I11").Value = 210.58333
$ws.Range("K11").Value = 210.58333
$ws.Range("M11").Value = -70.58332999999999
$ws.Range("H62").Value = 1555
$ws.Range("I62").Value = 1555
$ws.Range("K62").Value = 1555
$ws.Range("M62").Value = -931
$ws.Range("H65").Value = 1555
$ws.Range("I65").Value = 1555
$ws.Range("K65").Value = 7775
$ws.Range("M65").Value = -4655
$ws.Range("H92").Value = 1225.4375
$ws.Range("I92").Value = 1225.4375
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1225.4375
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 22.5625
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 1675.4166
$ws.Range("I96").Value = 1463.25
$ws.Range("K96").Value = 4389.75
$ws.Range("M96").Value = -3016.75
$ws.Range("H100").Value = 1392.9615
$ws.Range("J100").Value = 765.3333
$ws.Range("L100").Value = 765.3333
$ws.Range("N100").Value = -1847.3333
$ws.Range("H101").Value = 233.33333
$ws.Range("J101").Value = 300
$ws.Range("L101").Value = 900
$ws.Range("N101").Value = -4144
$ws.Range("H109").Value = 47500
$ws.Range("J109").Value = 47500
$ws.Range("L109").Value = 47500
$ws.Range("N109").Value = -50274
$ws.Range("H135").Value = 914.71875
$ws.Range("I135").Value = 784.7778
$ws.Range("K135").Value = 7063.000199999999
$ws.Range("M135").Value = -4528.000199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 53.1
$ws.Range("I5").Value = 40.142857
$ws.Range("J5").Value = 83.333336
$ws.Range("K5").Value = 40.142857
$ws.Range("L5").Value = 83.333336
$ws.Range("M5").Value = 71.85714300000001
$ws.Range("N5").Value = -307.333336
$ws.Range("H32").Value = 3105.123
$ws.Range("I32").Value = 3013.0156
$ws.Range("K32").Value = 3013.0156
$ws.Range("M32").Value = -2726.0156
$ws.Range("H97").Value = 2364.5
$ws.Range("I97").Value = 1712.7646
$ws.Range("K97").Value = 1712.7646
$ws.Range("M97").Value = -1216.7646
$ws.Range("H102").Value = 1603.3572
$ws.Range("I102").Value = 1378.0834
$ws.Range("K102").Value = 1378.0834
$ws.Range("M102").Value = 243.9166
$ws.Range("H122").Value = 2248.1304
$ws.Range("I122").Value = 2299.6316
$ws.Range("J122").Value = 2003.5
$ws.Range("K122").Value = 6898.8948
$ws.Range("L122").Value = 6010.5
$ws.Range("M122").Value = -4448.8948
$ws.Range("N122").Value = -10910.5
$ws.Range("H132").Value = 4142.3706
$ws.Range("I132").Value = 3918.3333
$ws.Range("K132").Value = 11754.9999
$ws.Range("M132").Value = -9224.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 53.1
$ws.Range("I4").Value = 40.142857
$ws.Range("J4").Value = 83.333336
$ws.Range("K4").Value = 40.142857
$ws.Range("L4").Value = 83.333336
$ws.Range("M4").Value = 74.85714300000001
$ws.Range("N4").Value = -313.333336
$ws.Range("H94").Value = 996.3333
$ws.Range("I94").Value = 522.2632
$ws.Range("K94").Value = 522.2632
$ws.Range("M94").Value = -71.26319999999998
$ws.Range("H99").Value = 3299.2173
$ws.Range("I99").Value = 1708.6666
$ws.Range("K99").Value = 1708.6666
$ws.Range("M99").Value = -210.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 982
$ws.Range("I16").Value = 917.7368
$ws.Range("J16").Value = 1389
$ws.Range("K16").Value = 917.7368
$ws.Range("L16").Value = 1389
$ws.Range("M16").Value = -630.7368
$ws.Range("N16").Value = -1963
$ws.Range("H99").Value = 3676.625
$ws.Range("I99").Value = 3120.2354
$ws.Range("J99").Value = 5027.857
$ws.Range("K99").Value = 3120.2354
$ws.Range("L99").Value = 5027.857
$ws.Range("M99").Value = -1622.2354
$ws.Range("N99").Value = -8023.857
$ws.Range("H105").Value = 1345.3334
$ws.Range("I105").Value = 1056.375
$ws.Range("J105").Value = 1923.25
$ws.Range("K105").Value = 1056.375
$ws.Range("L105").Value = 1923.25
$ws.Range("M105").Value = 690.625
$ws.Range("N105").Value = -5417.25
$ws.Range("H113").Value = 982
$ws.Range("I113").Value = 917.7368
$ws.Range("J113").Value = 1389
$ws.Range("K113").Value = 917.7368
$ws.Range("L113").Value = 1389
$ws.Range("M113").Value = 1252.2632
$ws.Range("N113").Value = -5729
$ws.Range("H126").Value = 3676.625
$ws.Range("I126").Value = 3120.2354
$ws.Range("J126").Value = 5027.857
$ws.Range("K126").Value = 9360.706200000001
$ws.Range("L126").Value = 15083.571
$ws.Range("M126").Value = -6890.706200000001
$ws.Range("N126").Value = -20023.571
$ws.Range("H132").Value = 2416.1428
$ws.Range("I132").Value = 2386.25
$ws.Range("K132").Value = 7158.75
$ws.Range("M132").Value = -4628.75
$ws.Range("H134").Value = 19509.611
$ws.Range("I134").Value = 5430.4287
$ws.Range("K134").Value = 16291.2861
$ws.Range("M134").Value = -13756.2861
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1459.2858
$ws.Range("J97").Value = 1508.1818
$ws.Range("L97").Value = 4524.5454
$ws.Range("N97").Value = -5516.5454
$ws.Range("H103").Value = 250.75
$ws.Range("I103").Value = 309.33334
$ws.Range("J103").Value = 75
$ws.Range("K103").Value = 928.0000200000001
$ws.Range("L103").Value = 225
$ws.Range("M103").Value = -49.00002000000006
$ws.Range("N103").Value = -1983
$ws.Range("H107").Value = 805.8889
$ws.Range("I107").Value = 621.5
$ws.Range("K107").Value = 1864.5
$ws.Range("M107").Value = 55.5
$ws.Range("H113").Value = 656.75
$ws.Range("I113").Value = 858.2
$ws.Range("J113").Value = 589.6
$ws.Range("K113").Value = 2574.6
$ws.Range("L113").Value = 1768.8
$ws.Range("M113").Value = -404.6000000000004
$ws.Range("N113").Value = -6108.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4005000
$ws.Range("I21").Value = 4005000
$ws.Range("K21").Value = 4005000
$ws.Range("M21").Value = -4004827
$ws.Range("H30").Value = 4005000
$ws.Range("I30").Value = 4005000
$ws.Range("K30").Value = 4005000
$ws.Range("M30").Value = -4004895
$ws.Range("H53").Value = 46000
$ws.Range("J53").Value = 46000
$ws.Range("L53").Value = 46000
$ws.Range("N53").Value = -47262
$ws.Range("H97").Value = 1023.44446
$ws.Range("I97").Value = 1175.6
$ws.Range("J97").Value = 833.25
$ws.Range("K97").Value = 1175.6
$ws.Range("L97").Value = 833.25
$ws.Range("M97").Value = -679.5999999999999
$ws.Range("N97").Value = -1825.25
$ws.Range("H102").Value = 26319952
$ws.Range("I102").Value = 1638.6571
$ws.Range("K102").Value = 1638.6571
$ws.Range("M102").Value = -16.6570999999999
$ws.Range("H122").Value = 4377.3
$ws.Range("I122").Value = 4162.8335
$ws.Range("K122").Value = 12488.5005
$ws.Range("M122").Value = -10038.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 175.26086
$ws.Range("I55").Value = 154.8125
$ws.Range("K55").Value = 154.8125
$ws.Range("M55").Value = 18.1875
$ws.Range("H93").Value = 1223
$ws.Range("I93").Value = 1188.7826
$ws.Range("J93").Value = 1380.4
$ws.Range("K93").Value = 1188.7826
$ws.Range("L93").Value = 1380.4
$ws.Range("M93").Value = 59.2174
$ws.Range("N93").Value = -3876.4
$ws.Range("H100").Value = 7277.4443
$ws.Range("I100").Value = 1666.3334
$ws.Range("J100").Value = 18499.666
$ws.Range("K100").Value = 1666.3334
$ws.Range("L100").Value = 18499.666
$ws.Range("M100").Value = -1125.3334
$ws.Range("N100").Value = -19581.666
$ws.Range("H122").Value = 281324.06
$ws.Range("I122").Value = 457130.03
$ws.Range("J122").Value = 5057.5
$ws.Range("K122").Value = 1371390.09
$ws.Range("L122").Value = 15172.5
$ws.Range("M122").Value = -1368940.09
$ws.Range("N122").Value = -20072.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 974.7143
$ws.Range("J96").Value = 947
$ws.Range("L96").Value = 947
$ws.Range("N96").Value = -3693
$ws.Range("H100").Value = 627.7143
$ws.Range("I100").Value = 666.25
$ws.Range("J100").Value = 396.5
$ws.Range("K100").Value = 1332.5
$ws.Range("L100").Value = 793
$ws.Range("M100").Value = -791.5
$ws.Range("N100").Value = -1875
$ws.Range("H132").Value = 2616.75
$ws.Range("I132").Value = 2497.1904
$ws.Range("K132").Value = 7491.5712
$ws.Range("M132").Value = -4961.5712

Write-Output "Applied all Coeurl_Profits market-data updates"